$d = $word.ActiveDocument

# 1) Update the report date
$d.Content.Find.Execute("23 Aug 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "26 Aug 2024", 2) | Out-Null

# 2) Clear the git revision hash (leaves the run, but empties its text)
$revRange = $d.Content
$revRange.Find.Execute("3aa81b4") | Out-Null
$revRange.Text = ""

# 3) Insert the new exploratory-analysis table after the paragraph that ends
#    "...and the school's existing ventilation system setting." and before
#    the "TODO: Add results." paragraph, followed by a blank paragraph.
$anchor = $d.Content
$anchor.Find.Execute("and the school's existing ventilation system setting.") | Out-Null
$targetPara = $anchor.Paragraphs(1)
$insertPoint = $d.Range($targetPara.Range.End, $targetPara.Range.End)

$tableXml = @"
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="TableGrid"/>
        <w:tblW w:w="0" w:type="auto"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
        <w:gridCol w:w="780"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:gridSpan w:val="12"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Time-varying covariate adjusted estimates of treatment effect for the primary outcome</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">None</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Ceiling</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Portable</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:gridSpan w:val="3"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Ceiling</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:gridSpan w:val="3"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="center"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Portable</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Superiority of air purification</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Noninferiority (portable vs ceiling)</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Mean</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Mean</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Mean</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">RR</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">[95%</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">CI]</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">RR</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">[95%</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">CI]</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="nil" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="left"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">PM2.5</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">1.03</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">0.86</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">0.55</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">0.94</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">[0.79</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">1.12]</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">0.70</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">[0.64</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">0.76]</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">&lt;0.001</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="0" w:type="auto"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="single" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
            <w:vAlign w:val="bottom"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
              <w:jc w:val="right"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">&lt;0.001</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:gridSpan w:val="12"/>
            <w:tcBorders>
              <w:top w:val="single" w:color="000000"/>
              <w:left w:val="nil" w:color="000000"/>
              <w:bottom w:val="nil" w:color="000000"/>
              <w:right w:val="nil" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:i w:val="true"/>
                <w:sz w:val="18"/>
              </w:rPr>
              <w:t xml:space="preserve">Sample means are unadjusted and do no account for the crossover design. Rate ratios (RRs) are adjusted for the crossover design, sensor type, first-order autocorrelation, the time-varying covariates selected by the cross-validated lasso, and clustering within classroom. RR &lt; 1 disfavors the reference (no air purification).</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p>
      <w:r>
        <w:t/>
      </w:r>
    </w:p>
</pkg:xmlData>
"@

$insertPoint.InsertXML($tableXml)
